$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-22 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-23 Monday", 2) | Out-Null
$d.Content.Find.Execute("643×6=3858", $true, $false, $false, $false, $false, $true, 1, $false, "929×9=8361", 2) | Out-Null
$d.Content.Find.Execute("425×7=2975", $true, $false, $false, $false, $false, $true, 1, $false, "497×4=1988", 2) | Out-Null
$d.Content.Find.Execute("116×2=232", $true, $false, $false, $false, $false, $true, 1, $false, "601×7=4207", 2) | Out-Null
$d.Content.Find.Execute("779×2=1558", $true, $false, $false, $false, $false, $true, 1, $false, "709×8=5672", 2) | Out-Null
$d.Content.Find.Execute("937×7=6559", $true, $false, $false, $false, $false, $true, 1, $false, "753×4=3012", 2) | Out-Null
$d.Content.Find.Execute("717×4=2868", $true, $false, $false, $false, $false, $true, 1, $false, "901×6=5406", 2) | Out-Null
$d.Content.Find.Execute("570×2=1140", $true, $false, $false, $false, $false, $true, 1, $false, "608×5=3040", 2) | Out-Null
$d.Content.Find.Execute("417×6=2502", $true, $false, $false, $false, $false, $true, 1, $false, "994×2=1988", 2) | Out-Null
$d.Content.Find.Execute("136×9=1224", $true, $false, $false, $false, $false, $true, 1, $false, "912×3=2736", 2) | Out-Null
$d.Content.Find.Execute("360×8=2880", $true, $false, $false, $false, $false, $true, 1, $false, "148×2=296", 2) | Out-Null
$d.Content.Find.Execute("659×4=2636", $true, $false, $false, $false, $false, $true, 1, $false, "208×8=1664", 2) | Out-Null
$d.Content.Find.Execute("930×2=1860", $true, $false, $false, $false, $false, $true, 1, $false, "651×4=2604", 2) | Out-Null
$d.Content.Find.Execute("199×7=1393", $true, $false, $false, $false, $false, $true, 1, $false, "424×6=2544", 2) | Out-Null
$d.Content.Find.Execute("545×6=3270", $true, $false, $false, $false, $false, $true, 1, $false, "597×6=3582", 2) | Out-Null
$d.Content.Find.Execute("658×9=5922", $true, $false, $false, $false, $false, $true, 1, $false, "626×9=5634", 2) | Out-Null
$d.Content.Find.Execute("929×4=3716", $true, $false, $false, $false, $false, $true, 1, $false, "457×7=3199", 2) | Out-Null
$d.Content.Find.Execute("435×8=3480", $true, $false, $false, $false, $false, $true, 1, $false, "529×2=1058", 2) | Out-Null
$d.Content.Find.Execute("214×8=1712", $true, $false, $false, $false, $false, $true, 1, $false, "397×5=1985", 2) | Out-Null
$d.Content.Find.Execute("867×6=5202", $true, $false, $false, $false, $false, $true, 1, $false, "558×4=2232", 2) | Out-Null
$d.Content.Find.Execute("339×9=3051", $true, $false, $false, $false, $false, $true, 1, $false, "441×7=3087", 2) | Out-Null
$d.Content.Find.Execute("428×6=2568", $true, $false, $false, $false, $false, $true, 1, $false, "874×9=7866", 2) | Out-Null
$d.Content.Find.Execute("518×7=3626", $true, $false, $false, $false, $false, $true, 1, $false, "582×6=3492", 2) | Out-Null
$d.Content.Find.Execute("908×2=1816", $true, $false, $false, $false, $false, $true, 1, $false, "962×7=6734", 2) | Out-Null
$d.Content.Find.Execute("846×6=5076", $true, $false, $false, $false, $false, $true, 1, $false, "800×3=2400", 2) | Out-Null
$d.Content.Find.Execute("458×9=4122", $true, $false, $false, $false, $false, $true, 1, $false, "525×4=2100", 2) | Out-Null
